$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(34, 34, 34, 34, 22, 6, 767, 787, 677)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("A13").Formula = "=SUM(A4:A12)"

$ws.Range("A14").Select()
